$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    279,280,
    313,314,
    317,318,
    346,347,
    351,352,
    372,373,
    379,380,
    382,383,
    389,390,
    400,401,
    419,420,
    421,422,
    431,432,
    457,458,
    536,537,
    579,580,
    583,584,
    586,587,
    593,594,
    601,602,
    720,721
)

for ($i = 0; $i -lt $pairs.Length; $i += 2) {
    $r1 = $pairs[$i]
    $r2 = $pairs[$i + 1]

    $range1 = $ws.Range("B${r1}:G${r1}")
    $range2 = $ws.Range("B${r2}:G${r2}")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value2 = $v2
    $range2.Value2 = $v1
}
